# Insert a new data row at row 4 (pushing existing rows 4..47 down to 5..48)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new record.
# Non price/date fields mirror the (now shifted) row directly below it,
# which held this record's original sibling values.
$ws.Cells.Item(4, 1).Value = 7
$ws.Cells.Item(4, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(4, 3).Value = "Ñuble"
$ws.Cells.Item(4, 4).Value = 44545
$ws.Cells.Item(4, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(4, 5).Value = 16
$ws.Cells.Item(4, 6).Value = 100112031
$ws.Cells.Item(4, 7).Value = "Poroto verde"
$ws.Cells.Item(4, 8).Value = "Sin especificar"
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 160
$ws.Cells.Item(4, 11).Value = 12500
$ws.Cells.Item(4, 12).Value = 13000
$ws.Cells.Item(4, 13).Value = 12750
$ws.Cells.Item(4, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(4, 15).Value = "Región del Maule"
$ws.Cells.Item(4, 16).Value = 510
$ws.Cells.Item(4, 17).Value = 25
$ws.Cells.Item(4, 18).Value = "Hortaliza"
